$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7, shifting existing rows 7-11 down to 8-12.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the weekly record (2021-12-03) for
# "Feria Lagunitas de Puerto Montt" / Arándano (blue), same as the
# surrounding rows for the columns that don't vary.
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C7").Value = "Los Lagos"
$ws.Range("D7").Value = 44533
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100101
$ws.Range("H7").Value = "Berries"
$ws.Range("I7").Value = 100101001
$ws.Range("J7").Value = "Arándano (blue)"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 400
$ws.Range("N7").Value = 3500
$ws.Range("O7").Value = 3600
$ws.Range("P7").Value = 3550
$ws.Range("Q7").Value = "`$/kilo"
$ws.Range("R7").Value = "Región del Maule"
$ws.Range("S7").Value = 3550
$ws.Range("T7").Value = 1
